# Rename all 30 worksheets to their new "summ<number>" names, preserving
# order/position, sheetId, and relationship id (only the <sheet name="..."/>
# attribute changes per the diff).

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ26261670",
    "summ26397965",
    "summ26534992",
    "summ26680974",
    "summ26815046",
    "summ26951541",
    "summ27095822",
    "summ27235853",
    "summ27382378",
    "summ27545844",
    "summ27784110",
    "summ27968882",
    "summ28167363",
    "summ28350775",
    "summ28517639",
    "summ28664410",
    "summ28811571",
    "summ28951080",
    "summ29083595",
    "summ29327198",
    "summ29466714",
    "summ29604233",
    "summ29744457",
    "summ29882641",
    "summ30029672",
    "summ30168073",
    "summ30309097",
    "summ30448924",
    "summ30594444",
    "summ30732465"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
